$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header values in row 1
$ws.Range("A1").Value = "relattion"
$ws.Range("B1").Value = "count"

# Set column A width to match the new custom width (stored OOXML "width" is
# ColumnWidth + 5/6 in this runtime's unit conversion, so back that off here
# to land exactly on width="36.5" in the saved file).
$ws.Columns.Item(1).ColumnWidth = 35.666666666666664

# Update the active selection to D7 (as recorded in the saved view state)
$ws.Range("D7").Select()
